# Update generated output data (bilibili manga/comic convention info)
# on the "展览" and "全部类型" worksheets to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 476
$ws1.Range("F3").Value = 5711
$ws1.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202409/naKbsgO31727426722022.jpeg"
$ws1.Range("F8").Value = 55
$ws1.Range("F9").Value = 543

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 476
$ws4.Range("F3").Value = 5711
$ws4.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202409/naKbsgO31727426722022.jpeg"
$ws4.Range("F10").Value = 55
$ws4.Range("F11").Value = 543
